$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '61.102.87'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -0.32%  '

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.322.49'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -0.52%  '

$ws.Range('E4').Value = '  -0.09%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '396.18'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -3.73%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '124.19'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +7.12%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.583'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +1.65%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.999'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -0.01%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.646'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +2.72%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.117'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +1.74%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '40.14'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +0.48%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '3.834.23'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -0.80%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '8.12'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -2.48%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '18.98'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -1.03%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '3.355.87'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +0.78%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '60.943.58'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -0.21%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.990'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -1.74%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '10.95'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +1.09%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.0000122'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +4.41%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '3.15'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -5.99%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '78.73'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +5.88%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '12.51'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -0.09%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '294.60'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +0.23%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '3.12'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +0.02%  '

$ws.Range('E26').Value = '  +9.87%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '28.53'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -2.03%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '8.02'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +6.00%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '7.33'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -7.08%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.170'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -1.37%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.00'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +0.04%  '

$ws.Range('E32').Value = '  -2.55%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '11.09'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -2.94%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '2.48'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -1.73%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '40.61'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -4.28%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.0469'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -4.30%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '51.74'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -1.13%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.995'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -0.36%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '3.32'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -3.76%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.84'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -8.30%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '135.35'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +0.85%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.93'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +1.42%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.120'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -0.26%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.274'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -4.76%  '

$ws.Range('B45').Value = 'NEARProtocol'
$ws.Range('C45').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '3.78'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -2.10%  '

$ws.Range('B46').Value = 'Celestia'
$ws.Range('C46').Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '16.30'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -0.83%  '

$ws.Range('E47').Value = '  -0.84%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '20.90'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -1.19%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '3.651.08'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -0.52%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '2.080.18'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -3.60%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '2.27'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -5.25%  '
